# Actualización automática 2025-09-18 09:12:30
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M3").Value = 1672.61
$ws1.Range("D13").Value = 915.84
$ws1.Range("M13").Value = 9364.549999999999
$ws1.Range("E30").Value = 69.13
$ws1.Range("D60").Value = "2 de 58"
$ws1.Range("E60").Value = "3 de 58"
$ws1.Range("M60").Value = "6 de 58"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F3").Value = 1672.61
$ws2.Range("F13").Value = 10280.39
$ws2.Range("F30").Value = 69.13
$ws2.Range("F60").Value = 28105.66

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 1866.24
$ws3.Range("E3").Value = 3638.37890386263
$ws3.Range("F3").Value = 0.3390316446231085

$ws3.Range("D4").Value = 2041.37
$ws3.Range("E4").Value = -737.3413934184
$ws3.Range("F4").Value = 1.565433449616782

$ws3.Range("D12").Value = 16147.88
$ws3.Range("E12").Value = 16256.92
$ws3.Range("F12").Value = 0.4983175332049573

$ws3.Range("D15").Value = 28105.66
$ws3.Range("E15").Value = 22378.10705102521
$ws3.Range("F15").Value = 0.556726679520427
